# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 132 (pushing the existing
# rows 132-167 down to 133-168) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("132").Insert()

$ws.Range("A132").Value = 6
$ws.Range("B132").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C132").Value = "Metropolitana"
$ws.Range("D132").Value = 44588
$ws.Range("E132").Value = 13
$ws.Range("F132").Value = 100112001
$ws.Range("G132").Value = "Berenjena"
$ws.Range("H132").Value = "Sin especificar"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 400
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 10000
$ws.Range("M132").Value = 9575
$ws.Range("N132").Value = "$/caja 50 unidades"
$ws.Range("O132").Value = "Región de Arica y Parinacota"
$ws.Range("P132").Value = 192
$ws.Range("Q132").Value = 50
$ws.Range("R132").Value = "Hortaliza"
